$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1099MISCdata")

# Test case 5: new recipient/payer test-run identifiers (row 3, 4, 9)
$ws.Range("B3").Value = "Test0312202191417"
$ws.Range("B4").Value = "Test0312202190806"
$ws.Range("B9").Value = "74-2042568"

# Confirmation number (row 5) must stay text even though it looks numeric
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "201204782"
$ws.Range("B5").Style = "Normal"

# Held amount / Payer EIN number (row 16) is a plain numeric value
$ws.Range("B16").Value = 901864257
